$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (D) and "Volume(1h)" (E) figures refreshed by the
# GitHub Actions cryptos-list updater; rows 15/16 (ShibaInu <-> Litecoin)
# also swapped position in the ranking.
$updates = [ordered]@{
    'D2' = '26.157.42'
    'E2' = '  -0.22%  '
    'D3' = '1.669.66'
    'E3' = '  -0.78%  '
    'D4' = '1.003'
    'E4' = '  -0.41%  '
    'D5' = '209.88'
    'E5' = '  -3.04%  '
    'D6' = '0.5228'
    'E6' = '  -0.20%  '
    'D7' = '1.003'
    'E7' = '  -0.45%  '
    'D8' = '0.2624'
    'E8' = '  -2.62%  '
    'D9' = '0.06330'
    'E9' = '  -0.60%  '
    'D10' = '21.21'
    'E10' = '  -1.24%  '
    'D11' = '0.07546'
    'E11' = '  -1.00%  '
    'D12' = '1.674.93'
    'E12' = '  -0.87%  '
    'D13' = '4.445'
    'E13' = '  -1.51%  '
    'D14' = '0.5506'
    'E14' = '  -4.32%  '
    'B15' = 'Litecoin'
    'C15' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'D15' = '66.59'
    'E15' = '  +0.85%  '
    'B16' = 'ShibaInu'
    'C16' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D16' = '0.000008004'
    'E16' = '  -3.87%  '
    'D17' = '26.170.40'
    'E18' = '  -0.45%  '
    'D19' = '4.759'
    'E19' = '  -2.20%  '
    'D20' = '187.07'
    'E20' = '  -1.24%  '
    'D21' = '10.34'
    'E21' = '  -4.33%  '
    'D22' = '6.205'
    'E22' = '  -0.46%  '
    'E23' = '  -0.41%  '
    'D24' = '149.68'
    'E24' = '  +0.55%  '
    'E25' = '  -0.90%  '
    'D26' = '7.506'
    'E26' = '  -3.61%  '
    'E27' = '  +0.88%  '
    'D28' = '0.06350'
    'E28' = '  +1.03%  '
    'D29' = '1.351'
    'E29' = '  -1.81%  '
    'E30' = '  -3.08%  '
    'D31' = '3.514'
    'E31' = '  -1.60%  '
    'D32' = '3.417'
    'E32' = '  -4.16%  '
    'D33' = '1.648'
    'E33' = '  -1.87%  '
    'D34' = '1.005'
    'E34' = '  -1.89%  '
    'D35' = '0.6035'
    'E35' = '  -1.39%  '
    'D36' = '2.409'
    'E36' = '  -0.48%  '
    'E37' = '  -0.04%  '
    'D38' = '6.144'
    'E38' = '  -0.79%  '
    'D39' = '1.110.88'
    'E39' = '  +1.19%  '
    'D40' = '0.01614'
    'D41' = '0.8657'
    'E41' = '  -2.41%  '
    'E42' = '  -0.69%  '
    'D43' = '100.34'
    'E43' = '  -0.18%  '
    'D44' = '1.824.54'
    'E44' = '  -0.48%  '
    'D45' = '0.00000000109'
    'E45' = '  -0.94%  '
    'D46' = '55.54'
    'E46' = '  -3.08%  '
    'D47' = '1.002'
    'E47' = '  -0.39%  '
    'D48' = '8.052'
    'E48' = '  -0.31%  '
    'E49' = '  -0.85%  '
    'D50' = '0.4240'
    'E50' = '  -1.05%  '
    'D51' = '5.934'
    'E51' = '  -1.46%  '
}

# Cells in the Price column can look like plain numbers (e.g. "1.003",
# "21.21"). Force them to be stored as text first so Excel does not
# silently reinterpret them as numeric/date values, which matches the
# original workbook where these are inline text strings.
foreach ($cellRef in $updates.Keys) {
    if ($cellRef.StartsWith("D")) {
        $ws.Range($cellRef).NumberFormat = "@"
    }
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
